$d = $word.ActiveDocument

# --- Update the date/title line ---
# Assign directly to the paragraph's Range.Text (without a trailing
# paragraph-mark character) so only the run text is replaced in place and
# no extra paragraph gets inserted.
$d.Paragraphs.Item(1).Range.Text = "2024-05-22 Wednesday"

# --- Update the division problems in the table ---
# Each problem lives in its own table cell. This runtime's Find/Replace
# operates over the whole story regardless of the Range it is invoked on,
# so instead we assign straight to Cell.Range.Text, which *is* properly
# scoped to the individual cell and preserves the existing run formatting.
# This also sidesteps any cross-talk between cells whose old/new values
# happen to collide (e.g. "203÷8=" -> "951÷6=" while a different cell
# already holds the literal text "951÷6=").
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; New = "676÷5=" },
    @{ Row = 1;  Col = 2; New = "110÷9=" },
    @{ Row = 1;  Col = 3; New = "319÷7=" },
    @{ Row = 1;  Col = 4; New = "102÷7=" },
    @{ Row = 1;  Col = 5; New = "223÷8=" },

    @{ Row = 5;  Col = 1; New = "893÷6=" },
    @{ Row = 5;  Col = 2; New = "786÷8=" },
    @{ Row = 5;  Col = 3; New = "649÷2=" },
    @{ Row = 5;  Col = 4; New = "414÷2=" },
    @{ Row = 5;  Col = 5; New = "490÷8=" },

    @{ Row = 9;  Col = 1; New = "951÷6=" },
    @{ Row = 9;  Col = 2; New = "778÷8=" },
    @{ Row = 9;  Col = 3; New = "584÷6=" },
    @{ Row = 9;  Col = 4; New = "547÷8=" },
    @{ Row = 9;  Col = 5; New = "314÷9=" },

    @{ Row = 13; Col = 1; New = "977÷5=" },
    @{ Row = 13; Col = 2; New = "713÷6=" },
    @{ Row = 13; Col = 3; New = "155÷9=" },
    @{ Row = 13; Col = 4; New = "741÷7=" },
    @{ Row = 13; Col = 5; New = "378÷4=" },

    @{ Row = 17; Col = 1; New = "717÷5=" },
    @{ Row = 17; Col = 2; New = "342÷2=" },
    @{ Row = 17; Col = 3; New = "625÷2=" },
    @{ Row = 17; Col = 4; New = "586÷8=" },
    @{ Row = 17; Col = 5; New = "251÷5=" }
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.New
}
